# #5: property boat&car done
# Expand the "汽車" (car) sheet (sheet3) from columns A:G to A:N, adding the
# property_category / category / date / legislator_name / legislator_id /
# source_file / index columns (same trailing columns used on the other
# property sheets), and rename column C's header from the literal engine
# number to the new "capacity" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車

# ---- Header row (row 1) ----
# B1 (name) and D1..G1 (owner/register_date/register_reason/acquire_value)
# already existed; only C1's label changes from the stray literal number to
# the new "capacity" string, and H1:N1 are brand-new header cells.
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"
$ws.Cells.Item(1,8).Value = "property_category"
$ws.Cells.Item(1,9).Value = "category"
$ws.Cells.Item(1,10).Value = "date"
$ws.Cells.Item(1,11).Value = "legislator_name"
$ws.Cells.Item(1,12).Value = "legislator_id"
$ws.Cells.Item(1,13).Value = "source_file"
$ws.Cells.Item(1,14).Value = "index"

# New header cells need the same bold/centered/bordered look as the rest of
# row 1 (style index 1 in styles.xml: bold font, thin box border, centered
# horizontal, top vertical).
for ($col = 8; $col -le 14; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# ---- Data row (row 2) ----
$ws.Cells.Item(2,2).Value = "LEXUSES330"
$ws.Cells.Item(2,5).Value = "94年06月23日"
$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"
# "2011-12-16" looks like a date to Excel's smart-entry parser, which would
# otherwise silently reinterpret it as a date serial number + date format.
# A leading apostrophe forces literal text (the apostrophe itself is not
# stored), matching the plain "date" string column used on every other
# sheet; then reset the style so no stray quote-prefix formatting lingers.
$ws.Cells.Item(2,10).Value = "'2011-12-16"
$ws.Cells.Item(2,10).Style = "Normal"
$ws.Cells.Item(2,11).Value = "管碧玲"
$ws.Cells.Item(2,12).Value = 1374
$ws.Cells.Item(2,13).Value = "tmp32301"
$ws.Cells.Item(2,14).Value = 44
